$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'242.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'23.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.414"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05900"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'3.440"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'6.524"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8088"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.9351"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Value = "'0.07393"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03301"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03068"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09348"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'3.857"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.001570"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04686"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.0005911"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.005979"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Value = "'0.004901"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'0.00006803"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'3.563"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'2.127"
$ws.Range("D24").Style = "Normal"
$ws.Range("D40").Value = "'0.03966"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.006179"
$ws.Range("D41").Style = "Normal"
$ws.Range("D44").Value = "'0.009507"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005214"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'0.6702"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.002383"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("D50").Style = "Normal"
